# Update column G ("K") values on the active worksheet.
# New values were regenerated (K instead of Strike#) per the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(6, 6, 4, 3, 7, 6, 5, 5, 8, 5, 3, 6, 1, 4, 8, 7, 3, 1, 4, 6, 2, 2, 3, 4, 3)

$row = 2
foreach ($val in $newValues) {
    $ws.Range("G$row").Value = $val
    $row++
}
